# Applies diff #5: insurance, claim, debt, investment done
# - Adds metadata columns (category/date/legislator_name/legislator_id/source_file/index, etc.)
#   to the 保險 (insurance, sheet index 5) and 債務 (debt, sheet index 6) sheets,
#   matching the pattern already used on the other sheets (土地/建物/汽車/存款).

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

function CopyStyle($srcCell, $dstCell) {
    $srcCell.Copy()
    $dstCell.PasteSpecial($xlPasteFormats)
}

# Helper: write a value that must stay plain text even when it looks like a
# date (e.g. "2012-04-24"), then restore the normal (non-text) number format
# by copying the style from a reference cell.
function SetTextValue($cell, $text, $styleSrcCell) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    CopyStyle $styleSrcCell $cell
}

# ---------------------------------------------------------------------------
# Sheet 5: 保險 (insurance)
# ---------------------------------------------------------------------------
$wsIns = $wb.Worksheets.Item(5)

# Header row 1: relabel B1/C1/D1/E1, and add F1:K1
$wsIns.Cells.Item(1,2).Value = "company"
$wsIns.Cells.Item(1,3).Value = "name"
$wsIns.Cells.Item(1,4).Value = "owner"
$wsIns.Cells.Item(1,5).Value = "property_category"

$wsIns.Cells.Item(1,6).Value = "category"
$wsIns.Cells.Item(1,7).Value = "date"
$wsIns.Cells.Item(1,8).Value = "legislator_name"
$wsIns.Cells.Item(1,9).Value = "legislator_id"
$wsIns.Cells.Item(1,10).Value = "source_file"
$wsIns.Cells.Item(1,11).Value = "index"

CopyStyle $wsIns.Cells.Item(1,2) $wsIns.Cells.Item(1,6)
CopyStyle $wsIns.Cells.Item(1,2) $wsIns.Cells.Item(1,7)
CopyStyle $wsIns.Cells.Item(1,2) $wsIns.Cells.Item(1,8)
CopyStyle $wsIns.Cells.Item(1,2) $wsIns.Cells.Item(1,9)
CopyStyle $wsIns.Cells.Item(1,2) $wsIns.Cells.Item(1,10)
CopyStyle $wsIns.Cells.Item(1,2) $wsIns.Cells.Item(1,11)

# Row 2 data: B2/C2/D2 unchanged, E2 relabeled, and add F2:K2
$wsIns.Cells.Item(2,5).Value = "insurance"

$wsIns.Cells.Item(2,6).Value = "normal"
$wsIns.Cells.Item(2,8).Value = "李俊俋"
$wsIns.Cells.Item(2,9).Value = 1738
$wsIns.Cells.Item(2,10).Value = "tmp16861"
$wsIns.Cells.Item(2,11).Value = 88

CopyStyle $wsIns.Cells.Item(2,2) $wsIns.Cells.Item(2,6)
SetTextValue $wsIns.Cells.Item(2,7) "2012-04-24" $wsIns.Cells.Item(2,2)
CopyStyle $wsIns.Cells.Item(2,2) $wsIns.Cells.Item(2,8)
CopyStyle $wsIns.Cells.Item(2,2) $wsIns.Cells.Item(2,9)
CopyStyle $wsIns.Cells.Item(2,2) $wsIns.Cells.Item(2,10)
CopyStyle $wsIns.Cells.Item(2,2) $wsIns.Cells.Item(2,11)

# ---------------------------------------------------------------------------
# Sheet 6: 債務 (debt)
# ---------------------------------------------------------------------------
$wsDebt = $wb.Worksheets.Item(6)

# Header row 1: relabel B1:G1, and add H1:N1
$wsDebt.Cells.Item(1,2).Value = "species"
$wsDebt.Cells.Item(1,3).Value = "debtor"
$wsDebt.Cells.Item(1,4).Value = "owner"
$wsDebt.Cells.Item(1,5).Value = "total"
$wsDebt.Cells.Item(1,6).Value = "register_date"
$wsDebt.Cells.Item(1,7).Value = "register_reason"

$wsDebt.Cells.Item(1,8).Value = "property_category"
$wsDebt.Cells.Item(1,9).Value = "category"
$wsDebt.Cells.Item(1,10).Value = "date"
$wsDebt.Cells.Item(1,11).Value = "legislator_name"
$wsDebt.Cells.Item(1,12).Value = "legislator_id"
$wsDebt.Cells.Item(1,13).Value = "source_file"
$wsDebt.Cells.Item(1,14).Value = "index"

CopyStyle $wsDebt.Cells.Item(1,2) $wsDebt.Cells.Item(1,8)
CopyStyle $wsDebt.Cells.Item(1,2) $wsDebt.Cells.Item(1,9)
CopyStyle $wsDebt.Cells.Item(1,2) $wsDebt.Cells.Item(1,10)
CopyStyle $wsDebt.Cells.Item(1,2) $wsDebt.Cells.Item(1,11)
CopyStyle $wsDebt.Cells.Item(1,2) $wsDebt.Cells.Item(1,12)
CopyStyle $wsDebt.Cells.Item(1,2) $wsDebt.Cells.Item(1,13)
CopyStyle $wsDebt.Cells.Item(1,2) $wsDebt.Cells.Item(1,14)

# Row 2 data: B2:G2 unchanged values, and add H2:N2
$wsDebt.Cells.Item(2,8).Value = "debt"
$wsDebt.Cells.Item(2,9).Value = "normal"
$wsDebt.Cells.Item(2,11).Value = "李俊俋"
$wsDebt.Cells.Item(2,12).Value = 1738
$wsDebt.Cells.Item(2,13).Value = "tmp16861"
$wsDebt.Cells.Item(2,14).Value = 98

CopyStyle $wsDebt.Cells.Item(2,2) $wsDebt.Cells.Item(2,8)
CopyStyle $wsDebt.Cells.Item(2,2) $wsDebt.Cells.Item(2,9)
SetTextValue $wsDebt.Cells.Item(2,10) "2012-04-24" $wsDebt.Cells.Item(2,2)
CopyStyle $wsDebt.Cells.Item(2,2) $wsDebt.Cells.Item(2,11)
CopyStyle $wsDebt.Cells.Item(2,2) $wsDebt.Cells.Item(2,12)
CopyStyle $wsDebt.Cells.Item(2,2) $wsDebt.Cells.Item(2,13)
CopyStyle $wsDebt.Cells.Item(2,2) $wsDebt.Cells.Item(2,14)

# Row 3 data: B3:G3 unchanged values, and add H3:N3
$wsDebt.Cells.Item(3,8).Value = "debt"
$wsDebt.Cells.Item(3,9).Value = "normal"
$wsDebt.Cells.Item(3,11).Value = "李俊俋"
$wsDebt.Cells.Item(3,12).Value = 1738
$wsDebt.Cells.Item(3,13).Value = "tmp16861"
$wsDebt.Cells.Item(3,14).Value = 99

CopyStyle $wsDebt.Cells.Item(3,2) $wsDebt.Cells.Item(3,8)
CopyStyle $wsDebt.Cells.Item(3,2) $wsDebt.Cells.Item(3,9)
SetTextValue $wsDebt.Cells.Item(3,10) "2012-04-24" $wsDebt.Cells.Item(3,2)
CopyStyle $wsDebt.Cells.Item(3,2) $wsDebt.Cells.Item(3,11)
CopyStyle $wsDebt.Cells.Item(3,2) $wsDebt.Cells.Item(3,12)
CopyStyle $wsDebt.Cells.Item(3,2) $wsDebt.Cells.Item(3,13)
CopyStyle $wsDebt.Cells.Item(3,2) $wsDebt.Cells.Item(3,14)
